$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.502.34'
$ws.Range("E2").Value = '  +1.05%  '
$ws.Range("D3").Value = '1.921.63'
$ws.Range("E3").Value = '  +1.34%  '
$ws.Range("E4").Value = '  +0.42%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '325.79'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.67%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.005'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.43%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4844'
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4085'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.64%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08174'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.27%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.024'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.15%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '23.80'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +5.73%  '
$ws.Range("D12").Value = '1.929.83'
$ws.Range("E12").Value = '  +1.74%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.048'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.35%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.246'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.88%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '91.63'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.12%  '
$ws.Range("B16").Value = 'TRON'
$ws.Range("C16").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.06780'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.67%  '
$ws.Range("B17").Value = 'BinanceUSD'
$ws.Range("C17").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.006'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.47%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001040'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.76%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.80'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.87%  '
$ws.Range("E20").Value = '  +0.29%  '
$ws.Range("D21").Value = '29.528.75'
$ws.Range("E21").Value = '  +1.06%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.636'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.57%  '
$ws.Range("E23").Value = '  +2.02%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.182'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.75%  '
$ws.Range("D25").Value = '2.158.21'
$ws.Range("E25").Value = '  +1.88%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.729'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +11.92%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '156.91'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.71%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.14'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.69%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.127'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.87%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '120.58'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.91%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.025'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.30%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09559'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.21%  '
$ws.Range("E33").Value = '  +3.71%  '
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.572'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.91%  '
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.396'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.16%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02288'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.04%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06147'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.26%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.184'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.16%  '
$ws.Range("B39").Value = 'TheSandbox'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5986'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.10%  '
$ws.Range("B40").Value = 'Aptos'
$ws.Range("C40").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '10.84'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +8.00%  '
$ws.Range("E41").Value = '  -0.55%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1863'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.08%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.438'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.32%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.281'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.00%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.07626'
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.49'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.74%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5592'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.25%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.964'
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '116.91'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.22%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.437'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +4.40%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.89'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.81%  '
